$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / Row 5 text + shared-string bookkeeping -------------------------
# Clear H6's text value (string date) first so that shared string slot is
# freed before touching J5/J6, keeping the resulting shared-string table
# ordering aligned with the target workbook.
$ws.Range("H6").ClearContents()

# Estatus (column J): row 6 goes from "Abierto" to "Proceso"
$ws.Range("J6").Value = "Proceso"

# Estatus (column J): row 5 goes from "Cerrado" to "Proceso"
$ws.Range("J5").Value = "Proceso"

# Acciones Correctivas (column E), row 6
$ws.Range("E6").Value = "Generar lo mas pronto posible las encuestas de satisfaccion"

# Desviaciones (column D), row 6
$ws.Range("D6").Value = "Por omicion a la generacion de encuestas de satisfaccion existen resultados invalidos en la seccion de satisfaccion"

# Responsable (column G), row 6
$ws.Range("G6").Value = "Jovanny Zepeda"

# --- Dates -------------------------------------------------------------
# Fecha de Deteccion (column F), row 6
$ws.Range("F6").Value = 42094

# Fecha Compromiso (column H), row 6 becomes a real date value (was text)
$ws.Range("H6").Value = 42124

# Fecha Real de Cierre (column I), row 5 cleared (no longer closed)
$ws.Range("I5").ClearContents()

# --- Formatting -----------------------------------------------------------
# Give D5 the same thin border already used throughout the table.
$b = $ws.Range("D5").Borders
$b.ColorIndex = 1
$b.LineStyle = 1

# Row 6 grows taller to fit the longer wrapped text.
$ws.Rows.Item(6).RowHeight = 120

# Move the active selection to J7.
$null = $ws.Range("J7").Select()

"done"
